# fix bug exeded requeste in google drive
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Bump the date in A1 by one day (45310 -> 45311)
$ws.Range("A1").Value = 45311

# Update price column (D) for rows 32-39
$ws.Range("D32").Value = 111.234
$ws.Range("D33").Value = 159.163
$ws.Range("D34").Value = 211.913
$ws.Range("D35").Value = 218.175
$ws.Range("D36").Value = 285.93
$ws.Range("D37").Value = 326.947
$ws.Range("D38").Value = 391.516
$ws.Range("D39").Value = 466.211
